$wb = $excel.ActiveWorkbook

$oldGuid = "951e4666-90f1-4658-ad72-79caf5d96d46"
$newGuid = "61781170-bded-4dda-96f6-272ae9349dca"

$oldZhXlf = "$oldGuid.8e7cc7606299b77deeff3337534125c3b24005ba.zh-cn.xlf"
$newZhXlf = "$newGuid.4853d14e551d25d59463264f76f7f80957f7b2fc.zh-cn.xlf"

$oldDeXlf = "$oldGuid.8e7cc7606299b77deeff3337534125c3b24005ba.de-de.xlf"
$newDeXlf = "$newGuid.4853d14e551d25d59463264f76f7f80957f7b2fc.de-de.xlf"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")

# A2: markdown file name referenced by multiple sheets (shared text)
$wsOverview.Range("A2").Value = "$newGuid.md"

# B2: path (cell value) and hyperlink display text both use this text
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = "e2e\$newGuid.md"
}

# G2: Latest HO Xliff Generate Date (this text is shared with de-de!H2)
$wsOverview.Range("G2").Value = "2016-08-28 04:56:37"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

# A2: markdown file name (same shared text as Overview!A2)
$wsZh.Range("A2").Value = "$newGuid.md"

foreach ($h in $wsZh.Hyperlinks) {
    $h.TextToDisplay = "$newGuid.md"
}

$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = "2016-08-28 04:56:33"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

# A2: markdown file name (same shared text as Overview!A2)
$wsDe.Range("A2").Value = "$newGuid.md"

foreach ($h in $wsDe.Hyperlinks) {
    $h.TextToDisplay = "$newGuid.md"
}

$wsDe.Range("G2").Value = $newDeXlf

# H2 originally shared the exact same text as Overview!G2 ("2016-08-28 04:56:19"),
# so it must follow the same new value to stay consistent with the shared string.
$wsDe.Range("H2").Value = "2016-08-28 04:56:37"
